$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42/43: Mantle and FirstDigitalUSD swap positions, with updated Price/Volume
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.779"
$ws.Range("E43").Value = "  -3.95%  "

# Price / Volume updates across the remaining rows
$ws.Range("D2").Value = "'62.879.12"
$ws.Range("E2").Value = "  -2.43%  "
$ws.Range("D3").Value = "'3.405.89"
$ws.Range("E3").Value = "  -3.09%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'576.30"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").Value = "'126.78"
$ws.Range("E6").Value = "  -5.78%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'3.407.02"
$ws.Range("E8").Value = "  -3.07%  "
$ws.Range("D9").Value = "'0.477"
$ws.Range("E9").Value = "  -2.39%  "
$ws.Range("D10").Value = "'7.40"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").Value = "'3.992.65"
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("D16").Value = "'3.408.42"
$ws.Range("E16").Value = "  -2.97%  "
$ws.Range("D17").Value = "'62.958.17"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").Value = "'24.99"
$ws.Range("E18").Value = "  -3.64%  "
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("D20").Value = "'5.70"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("D21").Value = "'13.21"
$ws.Range("E21").Value = "  -3.11%  "
$ws.Range("D22").Value = "'378.27"
$ws.Range("E22").Value = "  -4.21%  "
$ws.Range("D23").Value = "'0.560"
$ws.Range("E23").Value = "  -2.86%  "
$ws.Range("D24").Value = "'3.543.59"
$ws.Range("E24").Value = "  -3.03%  "
$ws.Range("D25").Value = "'72.55"
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -7.38%  "
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'7.01"
$ws.Range("E29").Value = "  -5.69%  "
$ws.Range("E30").Value = "  -4.27%  "
$ws.Range("D31").Value = "'7.91"
$ws.Range("E31").Value = "  -4.83%  "
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D35").Value = "'3.433.26"
$ws.Range("E35").Value = "  -3.08%  "
$ws.Range("D36").Value = "'22.88"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("D37").Value = "'5.34"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("D39").Value = "'164.46"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").Value = "'1.50"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("E41").Value = "  -3.73%  "
$ws.Range("D44").Value = "'41.67"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("E45").Value = "  -3.58%  "
$ws.Range("D46").Value = "'1.58"
$ws.Range("E46").Value = "  -5.62%  "
$ws.Range("D47").Value = "'22.98"
$ws.Range("E47").Value = "  -10.27%  "
$ws.Range("E48").Value = "  -7.34%  "
$ws.Range("D49").Value = "'6.69"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("D50").Value = "'2.258.70"
$ws.Range("E50").Value = "  -6.07%  "
$ws.Range("D51").Value = "'0.860"
$ws.Range("E51").Value = "  -4.50%  "
